$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (2007年 data) entirely; rows 3-6 shift up to 2-5
$ws.Rows("2:2").Delete()
